# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" sheets to match the latest scrape output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F
$updates = @{
    2  = 2783
    3  = 735
    5  = 6670
    6  = 1482
    9  = 34
    10 = 86
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
